# Add `owndat` (ADDR1, City, State, Zip Code) and `asmt` (VALASM1, VALASM2,
# VALASM3) fields to the IC PIN-level reference file output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns continue on from the existing last column (W = 23rd column).
# Row 1 holds the source-table name, row 2 holds the field name within
# that table - mirror the pattern already used for the other columns.

# Fill the whole header row (table names) first, then the whole field-name
# row, matching how the source workbook was authored.

# Row 1 - source table names
$ws.Range("X1").Value = "OWNDAT"
$ws.Range("Y1").Value = "OWNDAT"
$ws.Range("Z1").Value = "OWNDAT"
$ws.Range("AA1").Value = "OWNDAT"
$ws.Range("AB1").Value = "ASMT"
$ws.Range("AC1").Value = "ASMT"
$ws.Range("AD1").Value = "ASMT"

# Row 2 - field names within each table
$ws.Range("X2").Value = "ADDR1"
$ws.Range("Y2").Value = "City"
$ws.Range("Z2").Value = "State"
$ws.Range("AA2").Value = "Zip Code"
$ws.Range("AB2").Value = "VALASM1"
$ws.Range("AC2").Value = "VALASM2"
$ws.Range("AD2").Value = "VALASM3"

# Match the formatting already applied to the rest of the header columns
# (yellow fill + thin border) by copying an existing cell's format over.
$src = $ws.Range("W1:W2")
$dst = $ws.Range("X1:AD2")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
